# Update countries & provincias Spain
# Applies the data refresh captured in the diff:
#  - Swap the country names/data for two row pairs (Guinea-Bisau/Benin
#    and Polinesia Francesa/Letonia swapped places in the source list)
#  - Refresh a batch of per-country statistics
#  - Bump the "datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 03:08"

# Estados Unidos
$ws.Range("B4").Value = 7549299
$ws.Range("C4").Value = 51379
$ws.Range("D4").Value = 4775887
$ws.Range("E4").Value = 2559889
$ws.Range("G4").Value = 863
$ws.Range("H4").Value = 213523

# Peru
$ws.Range("B9").Value = 821564
$ws.Range("C9").Value = 3267
$ws.Range("D9").Value = 695645
$ws.Range("E9").Value = 93310
$ws.Range("G9").Value = 74
$ws.Range("H9").Value = 32609

# Argentina
$ws.Range("B11").Value = 779689
$ws.Range("C11").Value = 14687
$ws.Range("D11").Value = 614515
$ws.Range("E11").Value = 144575
$ws.Range("G11").Value = 311
$ws.Range("H11").Value = 20599

# Jordania
$ws.Range("D97").Value = 4929
$ws.Range("E97").Value = 8642

# Mauritania
$ws.Range("E114").Value = 186
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 162

# Trinidad yTobago
$ws.Range("B133").Value = 4629
$ws.Range("C133").Value = 59
$ws.Range("E133").Value = 1824

# Row 150/151: the source list now orders "Guinea-Bisau" before "Benin",
# so row 150 becomes Guinea-Bisau (refreshed numbers) and row 151 becomes
# the old Benin row (unchanged numbers).
$ws.Range("A150").Value = "Guinea-Bisau"
$ws.Range("B150").Value = 2362
$ws.Range("C150").Value = 38
$ws.Range("D150").Value = 1549
$ws.Range("E150").Value = 774
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 39

$ws.Range("A151").Value = "Benin"
$ws.Range("B151").Value = 2357
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 1973
$ws.Range("E151").Value = 343
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 41

# Row 158/159: likewise "Polinesia Francesa" now sorts before "Letonia".
$ws.Range("A158").Value = "Polinesia Francesa"
$ws.Range("B158").Value = 1964
$ws.Range("C158").Value = 112
$ws.Range("D158").Value = 1555
$ws.Range("E158").Value = 401
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 8

$ws.Range("A159").Value = "Letonia"
$ws.Range("B159").Value = 1945
$ws.Range("C159").Value = 77
$ws.Range("D159").Value = 1307
$ws.Range("E159").Value = 600
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 38

# Niger
$ws.Range("B167").Value = 1198
$ws.Range("C167").Value = 1
$ws.Range("D167").Value = 1115

# Taiwan
$ws.Range("B175").Value = 517
$ws.Range("C175").Value = 2
$ws.Range("E175").Value = 26

# Bermudas
$ws.Range("D192").Value = 169
$ws.Range("E192").Value = 3
